# feat: add 2022-Q3 data
#
# Before:  Sheet1 = "总计"      Sheet2 = "2021-Q2"   Sheet3 = "2021-Q1"
# After:   Sheet1 = "总计"      Sheet2 = "2022-Q3"   Sheet3 = "2021-Q2"   Sheet4 = "2021-Q1"
#
# A brand new quarterly sheet ("2022-Q3") is inserted right after the
# "总计" (totals) summary sheet, the totals sheet gets a new row describing
# it, and the previously-existing "2021-Q2"/"2021-Q1" sheets simply shift
# one position to the right (their contents are untouched).

$wb = $excel.ActiveWorkbook

$wsTotal = $wb.Worksheets.Item(1)   # 总计
$wsQ2    = $wb.Worksheets.Item(2)   # 2021-Q2 (must stay untouched)
$wsQ1    = $wb.Worksheets.Item(3)   # 2021-Q1 (must stay untouched)

# Helper: write a value into a cell as literal TEXT, even if it looks like a
# number (e.g. "004685" or "15.28"), without leaving any custom number
# format / style behind on the cell.
$blankCell = $wsQ2.Range("ZZ500")
function Set-TextValue {
    param($range, [string]$text)
    $range.Value = "'" + $text
    $blankCell.Copy()
    $range.PasteSpecial(-4122) # xlPasteFormats - re-clears any formatting quirk from the quote-prefix
}

# ---------------------------------------------------------------------
# 1) Duplicate the "2021-Q2" sheet to become the base of the new
#    "2022-Q3" sheet, placing it right after "总计" (i.e. before the
#    existing "2021-Q2" sheet). This keeps the original "2021-Q2" and
#    "2021-Q1" sheets completely byte-for-byte untouched.
# ---------------------------------------------------------------------
$wsQ2.Copy($wsQ2)
$wsQ3New = $wb.Worksheets.Item(2)
$wsQ3New.Name = "2022-Q3"

# ---------------------------------------------------------------------
# 2) Update the contents of the new "2022-Q3" sheet with the new fund.
# ---------------------------------------------------------------------
$wsQ3New.Range("D1").Value = "基金规模"

Set-TextValue $wsQ3New.Range("B2") "004685"
Set-TextValue $wsQ3New.Range("C2") "金元顺安元启灵活配置混合"
Set-TextValue $wsQ3New.Range("D2") "15.28"
Set-TextValue $wsQ3New.Range("E2") "77.14"
Set-TextValue $wsQ3New.Range("F2") "0.88"
Set-TextValue $wsQ3New.Range("G2") "0.1345"
$wsQ3New.Range("H2").Value = 3

# ---------------------------------------------------------------------
# 3) Update the "总计" (totals) sheet: the existing two rows shift down by
#    one quarter, and a brand new "2022-Q3" row is written in their place.
# ---------------------------------------------------------------------
$wsTotal.Range("A3").Copy($wsTotal.Range("A4"))
$wsTotal.Range("A4").Value = 2
$wsTotal.Range("B4").Value = "2021-Q1"
$wsTotal.Range("C4").Value = 1
$wsTotal.Range("D4").Value = 0.46

$wsTotal.Range("B3").Value = "2021-Q2"
$wsTotal.Range("D3").Value = 0.39

$wsTotal.Range("B2").Value = "2022-Q3"
$wsTotal.Range("D2").Value = 0.13

# ---------------------------------------------------------------------
# 4) Restore "2021-Q1" as the active sheet (it was the active sheet
#    before the edit; duplicating "2021-Q2" above made the new sheet
#    active instead).
# ---------------------------------------------------------------------
$wb.Worksheets.Item("2021-Q1").Activate()
